# ---------------------------------------------------------------------------
# Applies the "tentando realizar calibracao" commit to params.xlsx:
#   - workbook:  turn on iterative calculation with a max-change of 1E-4
#   - params (sheet1): D2 1 (was 2); add 4 new parameter rows (57-60); fix
#     selection
#   - levers (sheet2): collapse the lever table down to a single row (C.1-.0,5)
#     and widen column D; fix selection
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # params
$ws2 = $wb.Worksheets.Item(2)   # levers

# --- workbook level calc settings ------------------------------------------
$excel.Iteration  = $true
$excel.MaxChange  = 0.0001
$excel.MaxIterations = 100

# --- sheet "levers": drop to a single combination, widen column D ----------
$ws2.Rows("3:7").Delete()
$ws2.Range("D2").Value = 0.5

[void]$ws2.Columns("D").AutoFit()
$ws2.Columns("D").ColumnWidth = 20.736979166666668

# --- sheet "params": calibration tweak + new rows ---------------------------
$ws1.Range("D2").Value = 1

$ws1.Range("A57").Value = "aInitialSharePlayers"
$ws1.Range("C57").Value = 0.5
$ws1.Range("D57").Value = 0.5

$ws1.Range("A58").Value = "aInitialReorderShare"
$ws1.Range("C58").Value = 0.3
$ws1.Range("D58").Value = 0.4

$ws1.Range("A59").Value = "aTotalInitialInstalledBase"
$ws1.Range("C59").Value = 2500
$ws1.Range("D59").Value = 2500

$ws1.Range("A60").Value = "aInitialIndustryShipments"
$ws1.Range("C60").Value = 1831
$ws1.Range("D60").Value = 1831

# --- selections / active sheet, matching the saved view state --------------
$ws2.Activate()
[void]$ws2.Range("C6").Select()

$ws1.Activate()
[void]$ws1.Range("C53").Select()
